# Fix/normalize the "Recorded By" (column G) values on the active sheet.
# The original export sometimes listed "System" first; this pass re-orders
# the comma-separated recorder list so "System"/"system" is no longer
# always first (matching the upstream canonical export).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    switch ($val) {
        "System, admin@admin.com" {
            $cell.Value2 = "admin@admin.com, System"
        }
        "System, dnasr281@gmail.com" {
            $cell.Value2 = "dnasr281@gmail.com, System"
        }
        "System, system, backup@backdoor.com" {
            $cell.Value2 = "system, System, backup@backdoor.com"
        }
        "dnasr281@gmail.com, admin@admin.com" {
            $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
        }
        default { }
    }
}
